$wb = $excel.ActiveWorkbook

# --- "Prix Spot" sheet: add a new "27-jun" column (N) ---
$ws = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting from the previous header cell (M1) onto the new
# header cell (N1) so it picks up the same bold/border/centered style,
# then set its text.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("N1").Value = "27-jun"

$nValues = @(81.34, 66.17, 63.79, 45.14, 46.62, 71.3, 84.53, 90.09, 73.55, 37.04, 1.34, 0, -0.01, -0.01, -0.03, -0.01, 7.5, 20.06, 71.95, 104.61, 125.4, 125.49, 131.91, 102.48)

for ($i = 0; $i -lt $nValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 14).Value = $nValues[$i]
}

# --- "Gaz" sheet: update last price for 2025-06-25 ---
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("B8").Value = 34.75
